# "done with upload parcelle"
#
# The sheet originally held a single-farmer "parcelle" table with French
# headers (Nom, Details, Surface, plants, Plants Productive, Age Moyen,
# Production estimee, Production estimee vrac) plus an inspection marker.
#
# This edit reworks it into the uploaded/normalized parcelle shape used by
# the app: a new leading farmerId column, English/camelCase field names,
# and a second parcelle row ("Antohaka B") duplicating the first farmer's
# numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at A; everything that used to be in A..K shifts to B..L.
$ws.Columns.Item(1).Insert()

# Row 3 was a blank, pre-formatted spacer row. Clone row 2's formatting into
# it first so the new data we add lines up with row 2's cell styles exactly.
$ws.Range("A2:J2").Copy()
$ws.Range("A3:J3").PasteSpecial(-4122)

# Row 1: header labels (A..J)
$ws.Cells.Item(1,1).Value  = "farmerId"
$ws.Cells.Item(1,2).Value  = "name"
$ws.Cells.Item(1,3).Value  = "description"
$ws.Cells.Item(1,4).Value  = "overallSize"
$ws.Cells.Item(1,5).Value  = "totalPlants"
$ws.Cells.Item(1,6).Value  = "productivePlants"
$ws.Cells.Item(1,7).Value  = "averageAge"
$ws.Cells.Item(1,8).Value  = "estimatedProduction"
$ws.Cells.Item(1,9).Value  = "estimated_VRAC"
$ws.Cells.Item(1,10).Value = "inspected"

# Row 2: existing "Antohaka" parcelle, now prefixed with its farmerId
$ws.Cells.Item(2,1).Value  = 231
$ws.Cells.Item(2,2).Value  = "Antohaka"
$ws.Cells.Item(2,3).Value  = "Tegnarano"
$ws.Cells.Item(2,4).Value  = 0.97
$ws.Cells.Item(2,5).Value  = 2500
$ws.Cells.Item(2,6).Value  = 2000
$ws.Cells.Item(2,7).Value  = 6
$ws.Cells.Item(2,8).Value  = 700
$ws.Cells.Item(2,9).Value  = 175
$ws.Cells.Item(2,10).Value = $true

# Row 3: new "Antohaka B" parcelle for the same farmer
$ws.Cells.Item(3,1).Value  = 231
$ws.Cells.Item(3,2).Value  = "Antohaka B"
$ws.Cells.Item(3,3).Value  = "Tegnarano"
$ws.Cells.Item(3,4).Value  = 0.97
$ws.Cells.Item(3,5).Value  = 2500
$ws.Cells.Item(3,6).Value  = 2000
$ws.Cells.Item(3,7).Value  = 6
$ws.Cells.Item(3,8).Value  = 700
$ws.Cells.Item(3,9).Value  = 175
$ws.Cells.Item(3,10).Value = $true

# Row 4's spacer height shrank now that it no longer follows a taller block.
$ws.Rows.Item(4).RowHeight = 16.5
